# QA Round 2: deep quality optimization - compliance, diversification, UX improvements
# Co-authored-by: Cursor <cursoragent@cursor.com>

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Sheet restructuring: split the old "dickpic" sheet into two sheets --
#    a fresh, diversified "cumcontrol2" (content lives where "dickpic" used
#    to be) and a brand-new "dickpic" copy (placed right before "boosters"),
#    while the original "cumcontrol" sheet becomes "cumcontrol1".
# ---------------------------------------------------------------------------

$dickpic  = $wb.Worksheets.Item("dickpic")
$boosters = $wb.Worksheets.Item("boosters")

# Duplicate "dickpic" and drop the copy right before "boosters" so the tab
# order becomes: cumcontrol, dickpic, dickpic (2), boosters
$dickpic.Copy($boosters, $null)

# The original "dickpic" sheet is repurposed into "cumcontrol2"
$dickpic.Name = "cumcontrol2"

# The freshly-copied sheet becomes the new, untouched "dickpic"
$newDickpic = $wb.Worksheets.Item("dickpic (2)")
$newDickpic.Name = "dickpic"

# The original "cumcontrol" sheet becomes "cumcontrol1"
$cumcontrol1 = $wb.Worksheets.Item("cumcontrol")
$cumcontrol1.Name = "cumcontrol1"

# ---------------------------------------------------------------------------
# 2) VeraJourney copy tweaks
# ---------------------------------------------------------------------------

$journey = $wb.Worksheets.Item("VeraJourney")
$journey.Range("B10").Value = "hold on a sec"
$journey.Range("B11").Value = "I'm done holding back"
$journey.Range("B20").Value = "did you see it? 🥰"

# ---------------------------------------------------------------------------
# 3) "cumcontrol1" (previously "cumcontrol") copy tweaks
# ---------------------------------------------------------------------------

$cumcontrol1.Range("B2").Value = "if you finish before you see what I'm sending next you'll regret it"

$cumcontrol1.Range("B3").Value = "wait wait wait love... I have one more thing for you before you finish"
$cumcontrol1.Range("C3").Value = "DELAY. Send final PPV."

$cumcontrol1.Range("B4").Value = "I want to feel it at the same time... watch this first"
$cumcontrol1.Range("C4").Value = "SYNC variant. Send PPV."

$cumcontrol1.Range("B5").Value = "okay NOW we can go together... open this"
$cumcontrol1.Range("C5").Value = "SYNC. Send PPV."

$cumcontrol1.Range("B6").Value = "you better not be close already... I have more to show you"

$cumcontrol1.Range("B7").Value = "not yet... I said not yet love"
$cumcontrol1.Range("C7").Value = "CONTROL. More PPVs to send. Create urgency to open next."

# ---------------------------------------------------------------------------
# 4) "cumcontrol2" (previously "dickpic") full content rewrite
# ---------------------------------------------------------------------------

$cumcontrol2 = $wb.Worksheets.Item("cumcontrol2")

$cumcontrol2.Range("A2").Value = "delay2"
$cumcontrol2.Range("B2").Value = "hold on just a little longer love, I promise this next one is worth it"
$cumcontrol2.Range("C2").Value = "DELAY variant."

$cumcontrol2.Range("A3").Value = "delay1"
$cumcontrol2.Range("B3").Value = "don't you dare... not until you see what I just did"
$cumcontrol2.Range("C3").Value = "DELAY. Send PPV."

$cumcontrol2.Range("A4").Value = "sync2"
$cumcontrol2.Range("B4").Value = "let's do this together... but you have to open this first"
$cumcontrol2.Range("C4").Value = "SYNC variant."

$cumcontrol2.Range("A5").Value = "sync1"
$cumcontrol2.Range("B5").Value = "okay I'm ready now too love... watch this with me"
$cumcontrol2.Range("C5").Value = "SYNC. Send PPV."

$cumcontrol2.Range("A6").Value = "edge2"
$cumcontrol2.Range("B6").Value = "patience... the best part hasn't even happened yet"
$cumcontrol2.Range("C6").Value = "EDGE variant."

$cumcontrol2.Range("A7").Value = "edge1"
$cumcontrol2.Range("B7").Value = "slow down... I'm not letting you off that easy"
$cumcontrol2.Range("C7").Value = "CONTROL."
